$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the saved view: scroll position + selected cell ---
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L39").Select()

# --- Re-point the "small (%)" averages at the new pre-train run's numbers ---
$ws.Range("D29").Formula = "=(0.0409356725146198 + 0.0487329434697855 + 0.0487329434697855 + 0.0526315789473684 + 0.0448343079922027) / 5 * 100"
$ws.Range("D30").Formula = "=(0.0892966360856269 + 0.0951070336391437 + 0.0938837920489296 + 0.0889908256880734 + 0.0954128440366972) / 5 * 100"

$ws.Range("D32").Formula = "=(0.019047619047619 + 0.0380952380952381 + 0.0285714285714285 + 0.0285714285714285 + 0.019047619047619) / 5 * 100"
$ws.Range("D33").Formula = "=(0.0902255639097744 + 0.0796992481203007 + 0.0736842105263157 + 0.0721804511278195 + 0.069172932330827) / 5 * 100"

$ws.Range("D36").Formula = "=(0.0175438596491228 + 0.0584795321637426 + 0.0487329434697855 + 0.0253411306042884 + 0.0389863547758284) / 5 * 100"
$ws.Range("D37").Formula = "=(0.091131498470948 + 0.0877675840978593 + 0.0856269113149847 + 0.0877675840978593 + 0.0963302752293578) / 5 * 100"

$ws.Range("D39").Formula = "=(0.00952380952380952 + 0.0285714285714285 + 0.0380952380952381 + 0.0476190476190476 + 0.0571428571428571) / 5 * 100"
$ws.Range("D40").Formula = "=(0.0796992481203007 + 0.0857142857142857 + 0.0827067669172932 + 0.0947368421052631 + 0.0887218045112782) / 5 * 100"

$ws.Range("D43").Formula = "=(0.490712074303405 + 0.467492260061919 + 0.530959752321981 + 0.507739938080495 + 0.5) / 5 * 100"
$ws.Range("D44").Formula = "=(0.446163366336633 + 0.441212871287128 + 0.435024752475247 + 0.430074257425742 + 0.440903465346534) / 5 * 100"

$ws.Range("D46").Formula = "=(0.305343511450381 + 0.358778625954198 + 0.32824427480916 + 0.312977099236641 + 0.343511450381679) / 5 * 100"
$ws.Range("D47").Formula = "=(0.278963414634146 + 0.234756097560975 + 0.28810975609756 + 0.268292682926829 + 0.283536585365853) / 5 * 100"

$ws.Range("D50").Formula = "=(0.41640866873065 + 0.455108359133126 + 0.439628482972136 + 0.411764705882352 + 0.452012383900928) / 5 * 100"
$ws.Range("D51").Formula = "=(0.398205445544554 + 0.394183168316831 + 0.392017326732673 + 0.395111386138613 + 0.39789603960396) / 5 * 100"

$ws.Range("D53").Formula = "=(0.229007633587786 + 0.206106870229007 + 0.221374045801526 + 0.305343511450381 + 0.297709923664122) / 5 * 100"
$ws.Range("D54").Formula = "=(0.216463414634146 + 0.204268292682926 + 0.213414634146341 + 0.236280487804878 + 0.198170731707317) / 5 * 100"

$ws.Range("F55").Formula = "=(0.0993930197268588) / 1 * 100"
$ws.Range("F56").Formula = "=(0.11067476294194) / 1 * 100"
